# repull data, push all data, mean calculation
# Update column F (dSF) values for the rows whose scraped data changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    4  = -4
    6  = -2
    7  = 2
    10 = -2
    11 = -4
    12 = -3
    13 = 5
    15 = -2
    16 = 1
    17 = 9
    18 = -4
    19 = 1
    20 = 1
    22 = -1
    23 = -6
    24 = 3
    25 = -1
    26 = -1
    28 = 1
    29 = 1
    30 = 5
    31 = 1
    32 = 2
    34 = 6
    35 = 3
    36 = -5
    37 = 1
    38 = -1
    39 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
